# Auto-generated edit script: updates cached market-price snapshot values
# across the 8 job-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 507.5
$ws.Range("I18").Value = 507.5
$ws.Range("K18").Value = 507.5
$ws.Range("M18").Value = -223.5
$ws.Range("H74").Value = 8123.2666
$ws.Range("I74").Value = 7834.8
$ws.Range("J74").Value = 8700.200000000001
$ws.Range("K74").Value = 7834.8
$ws.Range("L74").Value = 8700.200000000001
$ws.Range("M74").Value = -6898.8
$ws.Range("N74").Value = -10572.2
$ws.Range("H76").Value = 7850
$ws.Range("I76").Value = 6950
$ws.Range("K76").Value = 6950
$ws.Range("M76").Value = -6635
$ws.Range("H77").Value = 8123.2666
$ws.Range("I77").Value = 7834.8
$ws.Range("J77").Value = 8700.200000000001
$ws.Range("K77").Value = 39174
$ws.Range("L77").Value = 43501
$ws.Range("M77").Value = -34494
$ws.Range("N77").Value = -52861
$ws.Range("H79").Value = 7850
$ws.Range("I79").Value = 6950
$ws.Range("K79").Value = 6950
$ws.Range("M79").Value = -5858
$ws.Range("H137").Value = 1214480.9
$ws.Range("I137").Value = 1033.3334
$ws.Range("K137").Value = 3100.0002
$ws.Range("M137").Value = -550.0001999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 46333
$ws.Range("J24").Value = 46333
$ws.Range("L24").Value = 46333
$ws.Range("N24").Value = -47081
$ws.Range("H61").Value = 3319.0952
$ws.Range("I61").Value = 2358.0908
$ws.Range("J61").Value = 4376.2
$ws.Range("K61").Value = 2358.0908
$ws.Range("L61").Value = 4376.2
$ws.Range("M61").Value = -2146.0908
$ws.Range("N61").Value = -4800.2
$ws.Range("H74").Value = 12237.941
$ws.Range("I74").Value = 1149.2858
$ws.Range("K74").Value = 1149.2858
$ws.Range("M74").Value = -275.2858000000001
$ws.Range("H77").Value = 12237.941
$ws.Range("I77").Value = 1149.2858
$ws.Range("K77").Value = 5746.429
$ws.Range("M77").Value = -1378.429
$ws.Range("H98").Value = 49500
$ws.Range("J98").Value = 49500
$ws.Range("L98").Value = 49500
$ws.Range("N98").Value = -55490
$ws.Range("H100").Value = 46333
$ws.Range("J100").Value = 46333
$ws.Range("L100").Value = 46333
$ws.Range("N100").Value = -48497
$ws.Range("H102").Value = 5413.407
$ws.Range("J102").Value = 4293.3335
$ws.Range("L102").Value = 4293.3335
$ws.Range("N102").Value = -7537.3335
$ws.Range("H122").Value = 3425.75
$ws.Range("I122").Value = 3645.0952
$ws.Range("J122").Value = 1890.3334
$ws.Range("K122").Value = 10935.2856
$ws.Range("L122").Value = 5671.0002
$ws.Range("M122").Value = -8485.285600000001
$ws.Range("N122").Value = -10571.0002
$ws.Range("H132").Value = 4852.8657
$ws.Range("I132").Value = 5191.434
$ws.Range("K132").Value = 15574.302
$ws.Range("M132").Value = -13044.302
$ws.Range("H136").Value = 3319.0952
$ws.Range("I136").Value = 2358.0908
$ws.Range("J136").Value = 4376.2
$ws.Range("K136").Value = 7074.2724
$ws.Range("L136").Value = 13128.6
$ws.Range("M136").Value = -4524.2724
$ws.Range("N136").Value = -18228.6

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 27175.666
$ws.Range("J95").Value = 27175.666
$ws.Range("L95").Value = 27175.666
$ws.Range("N95").Value = -32667.666
$ws.Range("H107").Value = 1591.7333
$ws.Range("J107").Value = 5000
$ws.Range("L107").Value = 5000
$ws.Range("N107").Value = -8840

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H122").Value = 1639.5
$ws.Range("I122").Value = 1683
$ws.Range("J122").Value = 1509
$ws.Range("K122").Value = 5049
$ws.Range("L122").Value = 4527
$ws.Range("M122").Value = -2599
$ws.Range("N122").Value = -9427

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 931.381
$ws.Range("J107").Value = 1084.8334
$ws.Range("L107").Value = 3254.5002
$ws.Range("N107").Value = -7094.5002
$ws.Range("H124").Value = 7374.5
$ws.Range("H132").Value = 1817.5555
$ws.Range("J132").Value = 1982.25
$ws.Range("L132").Value = 17840.25
$ws.Range("N132").Value = -22900.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H70").Value = 39387.617
$ws.Range("I70").Value = 43571.625
$ws.Range("J70").Value = 25998.8
$ws.Range("K70").Value = 43571.625
$ws.Range("L70").Value = 25998.8
$ws.Range("M70").Value = -43301.625
$ws.Range("N70").Value = -26538.8
$ws.Range("H73").Value = 39387.617
$ws.Range("I73").Value = 43571.625
$ws.Range("J73").Value = 25998.8
$ws.Range("K73").Value = 43571.625
$ws.Range("L73").Value = 25998.8
$ws.Range("M73").Value = -42635.625
$ws.Range("N73").Value = -27870.8
$ws.Range("H92").Value = 40049.2
$ws.Range("I92").Value = 29999
$ws.Range("J92").Value = 42561.75
$ws.Range("K92").Value = 29999
$ws.Range("L92").Value = 42561.75
$ws.Range("M92").Value = -28127
$ws.Range("N92").Value = -46305.75
$ws.Range("H132").Value = 4741.6924
$ws.Range("I132").Value = 4146.3887
$ws.Range("K132").Value = 12439.1661
$ws.Range("M132").Value = -9909.166100000002

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5825.119
$ws.Range("I132").Value = 5343.8613
$ws.Range("K132").Value = 16031.5839
$ws.Range("M132").Value = -13501.5839

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982
$ws.Range("H122").Value = 3236.205
$ws.Range("I122").Value = 3567.3872
$ws.Range("J122").Value = 1952.875
$ws.Range("K122").Value = 10702.1616
$ws.Range("L122").Value = 5858.625
$ws.Range("M122").Value = -8252.161599999999
$ws.Range("N122").Value = -10758.625
$ws.Range("H126").Value = 41670316
$ws.Range("I126").Value = 3951.3125
$ws.Range("J126").Value = 125003040
$ws.Range("K126").Value = 11853.9375
$ws.Range("L126").Value = 375009120
$ws.Range("M126").Value = -9383.9375
$ws.Range("N126").Value = -375014060
$ws.Range("H132").Value = 2626.5454
$ws.Range("I132").Value = 1783.8334
$ws.Range("J132").Value = 6418.75
$ws.Range("K132").Value = 5351.5002
$ws.Range("L132").Value = 19256.25
$ws.Range("M132").Value = -2821.5002
$ws.Range("N132").Value = -24316.25
